$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-02-27 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-02-28 Wednesday", 2)

# Update the division-fact table. The problems live in rows 1, 5, 9, 13, 17
# (5 columns each); the intervening rows are blank spacer rows. Addressing
# cells directly by (row, column) avoids any ambiguity from duplicate
# "NN÷N=" text occurring more than once in the document.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "66÷6="
$t.Cell(1, 2).Range.Text  = "10÷7="
$t.Cell(1, 3).Range.Text  = "48÷5="
$t.Cell(1, 4).Range.Text  = "28÷9="
$t.Cell(1, 5).Range.Text  = "82÷4="

$t.Cell(5, 1).Range.Text  = "84÷7="
$t.Cell(5, 2).Range.Text  = "85÷9="
$t.Cell(5, 3).Range.Text  = "70÷7="
$t.Cell(5, 4).Range.Text  = "94÷4="
$t.Cell(5, 5).Range.Text  = "94÷3="

$t.Cell(9, 1).Range.Text  = "83÷3="
$t.Cell(9, 2).Range.Text  = "41÷9="
# Cell (9,3) "78÷9=" is unchanged by the diff.
$t.Cell(9, 4).Range.Text  = "40÷3="
$t.Cell(9, 5).Range.Text  = "48÷7="

$t.Cell(13, 1).Range.Text = "14÷3="
$t.Cell(13, 2).Range.Text = "30÷5="
$t.Cell(13, 3).Range.Text = "96÷9="
$t.Cell(13, 4).Range.Text = "31÷8="
$t.Cell(13, 5).Range.Text = "76÷3="

$t.Cell(17, 1).Range.Text = "11÷3="
$t.Cell(17, 2).Range.Text = "40÷8="
$t.Cell(17, 3).Range.Text = "76÷9="
$t.Cell(17, 4).Range.Text = "81÷8="
$t.Cell(17, 5).Range.Text = "24÷9="
